$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the data (columns B:L) up by one row for rows 12-15, pulling in the
# values that previously lived in rows 13-16 respectively. Column A (the
# "Collection N" labels) stays put.
$row13 = $ws.Range("B13:L13").Value2
$row14 = $ws.Range("B14:L14").Value2
$row15 = $ws.Range("B15:L15").Value2
$row16 = $ws.Range("B16:L16").Value2

$ws.Range("B12:L12").Value2 = $row13
$ws.Range("B13:L13").Value2 = $row14
$ws.Range("B14:L14").Value2 = $row15
$ws.Range("B15:L15").Value2 = $row16

# The last row no longer has any data; clear it out (including the
# "Collection 19" label), leaving only the formatted-but-empty A16 cell.
$ws.Range("B16:L16").Clear()
$ws.Range("A16").ClearContents()

# Column A is now auto best-fit sized.
$ws.Columns("A").AutoFit()

# Cursor/selection moved to A11 in the saved file.
$ws.Range("A11").Select()
